$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowRange {
    param($ws, $rowNum, $values)
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) { $arr[0,$i] = $values[$i] }
    $rng = $ws.Range("F" + $rowNum + ":V" + $rowNum)
    $rng.Value = $arr
}

# --- Swap/rotate match data between rows (F:V), keep A:E (index/meta) untouched ---
$row5 = @("Roskilde", 1, "FA 2000", 0, 1.58, "04/08/2023 22:16", 1.61, "05/08/2023 13:01", 4.16, "04/08/2023 22:16", 4.53, "05/08/2023 13:01", 4.83, "04/08/2023 22:16", 4.34, "05/08/2023 13:01", "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-frederiksberg-alliancen-2000/fLdM47Ag/")
Set-RowRange $ws 5 $row5
$row6 = @("Brabrand", 1, "Thisted FC", 2, 1.95, "04/08/2023 02:12", 2.47, "05/08/2023 12:25", 3.36, "04/08/2023 02:12", 3.4, "05/08/2023 12:04", 3.26, "04/08/2023 02:12", 2.69, "05/08/2023 12:25", "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-thisted-fc/bg3D6Tus/")
Set-RowRange $ws 6 $row6
$row20 = @("FA 2000", 2, "AB Copenhagen", 2, 3.11, "25/08/2023 11:42", 3.13, "25/08/2023 18:44", 3.56, "25/08/2023 11:42", 3.75, "25/08/2023 18:44", 2.06, "25/08/2023 11:42", 2.06, "25/08/2023 18:44", "https://www.betexplorer.com/football/denmark/2nd-division/frederiksberg-alliancen-2000-ab-copenhagen/lMQNCM8U/")
Set-RowRange $ws 20 $row20
$row21 = @("Hellerup", 2, "Nykobing", 3, 2.33, "25/08/2023 11:42", 2.41, "25/08/2023 18:44", 3.42, "25/08/2023 11:42", 3.73, "25/08/2023 17:25", 2.81, "25/08/2023 11:42", 2.57, "25/08/2023 18:44", "https://www.betexplorer.com/football/denmark/2nd-division/hellerup-nykobing/EqhUX5vt/")
Set-RowRange $ws 21 $row21
$row26 = @("Skive", 2, "FA 2000", 2, 2.44, "02/09/2023 10:43", 2.58, "02/09/2023 13:48", 3.61, "02/09/2023 10:43", 3.43, "02/09/2023 13:48", 2.56, "02/09/2023 10:43", 2.55, "02/09/2023 13:48", "https://www.betexplorer.com/football/denmark/2nd-division/skive-frederiksberg-alliancen-2000/K2AsTsgB/")
Set-RowRange $ws 26 $row26
$row27 = @("Brabrand", 1, "Hellerup", 2, 2.51, "01/09/2023 02:12", 2.37, "01/09/2023 23:41", 3.34, "01/09/2023 02:12", 3.56, "02/09/2023 12:03", 2.45, "01/09/2023 02:12", 2.68, "01/09/2023 23:41", "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-hellerup/EP9oS18H/")
Set-RowRange $ws 27 $row27
$row32 = @("FA 2000", 2, "Hellerup", 0, 2.18, "08/09/2023 08:12", 2.08, "08/09/2023 17:30", 3.44, "08/09/2023 08:12", 3.56, "08/09/2023 18:32", 2.95, "08/09/2023 08:12", 3.02, "08/09/2023 17:30", "https://www.betexplorer.com/football/denmark/2nd-division/frederiksberg-alliancen-2000-hellerup/0Gr4pqOo/")
Set-RowRange $ws 32 $row32
$row33 = @("Thisted FC", 0, "Skive", 0, 2.16, "07/09/2023 07:12", 2.41, "08/09/2023 18:19", 3.37, "07/09/2023 07:12", 3.39, "08/09/2023 18:56", 2.9, "07/09/2023 07:12", 2.77, "08/09/2023 18:56", "https://www.betexplorer.com/football/denmark/2nd-division/thisted-fc-skive/MVWCrNhb/")
Set-RowRange $ws 33 $row33
$row35 = @("Esbjerg", 3, "Nykobing", 4, 1.36, "10/09/2023 09:12", 1.32, "10/09/2023 13:24", 5.13, "10/09/2023 09:12", 5.56, "10/09/2023 13:24", 6.23, "10/09/2023 09:12", 7.4, "10/09/2023 13:24", "https://www.betexplorer.com/football/denmark/2nd-division/esbjerg-nykobing/EwkKt1NA/")
Set-RowRange $ws 35 $row35
$row36 = @("Aarhus Fremad", 3, "F. Amager", 1, 1.34, "09/09/2023 02:12", 1.49, "10/09/2023 13:31", 4.74, "09/09/2023 02:12", 4.45, "10/09/2023 13:31", 6.73, "09/09/2023 02:12", 5.68, "10/09/2023 13:31", "https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-fremad-amager/hUmOuLxH/")
Set-RowRange $ws 36 $row36
$row44 = @("Aarhus Fremad", 3, "FA 2000", 0, 1.49, "22/09/2023 01:12", 1.5, "23/09/2023 13:45", 4.35, "22/09/2023 01:12", 4.56, "23/09/2023 13:45", 4.91, "22/09/2023 01:12", 5.43, "23/09/2023 13:45", "https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-frederiksberg-alliancen-2000/rHoerEbF/")
Set-RowRange $ws 44 $row44
$row45 = @("Middelfart", 1, "AB Copenhagen", 0, 1.95, "22/09/2023 01:12", 2.31, "23/09/2023 13:46", 3.49, "22/09/2023 01:12", 3.25, "23/09/2023 13:49", 3.17, "22/09/2023 01:12", 3.02, "23/09/2023 13:49", "https://www.betexplorer.com/football/denmark/2nd-division/middelfart-ab-copenhagen/M7pasYDL/")
Set-RowRange $ws 45 $row45
$row55 = @("Aarhus Fremad", 1, "AB Copenhagen", 1, 1.55, "06/10/2023 01:13", 1.48, "07/10/2023 13:58", 4.09, "06/10/2023 01:13", 4.81, "07/10/2023 13:58", 4.41, "06/10/2023 01:13", 5.36, "07/10/2023 13:58", "https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-ab-copenhagen/AXvIm9ze/")
Set-RowRange $ws 55 $row55
$row56 = @("Middelfart", 0, "Roskilde", 1, 2.18, "06/10/2023 01:12", 2.26, "07/10/2023 10:07", 3.21, "06/10/2023 01:12", 3.32, "07/10/2023 12:01", 2.88, "06/10/2023 01:12", 3.01, "07/10/2023 10:07", "https://www.betexplorer.com/football/denmark/2nd-division/middelfart-roskilde/S8kDlkLl/")
Set-RowRange $ws 56 $row56
$row57 = @("Skive", 0, "Nykobing", 2, 2.58, "06/10/2023 01:13", 3.55, "07/10/2023 13:40", 3.24, "06/10/2023 01:13", 3.55, "07/10/2023 13:40", 2.39, "06/10/2023 01:13", 1.97, "07/10/2023 13:40", "https://www.betexplorer.com/football/denmark/2nd-division/skive-nykobing/lIj9kV5r/")
Set-RowRange $ws 57 $row57
$row65 = @("Nykobing", 0, "Brabrand", 0, 1.56, "13/10/2023 02:13", 1.68, "14/10/2023 12:36", 3.95, "13/10/2023 02:13", 3.75, "14/10/2023 13:04", 4.5, "13/10/2023 02:13", 4.68, "14/10/2023 12:36", "https://www.betexplorer.com/football/denmark/2nd-division/nykobing-brabrand/ABwcw7Sg/")
Set-RowRange $ws 65 $row65
$row66 = @("AB Copenhagen", 0, "Skive", 0, 1.72, "13/10/2023 02:13", 1.62, "14/10/2023 14:51", 3.63, "13/10/2023 02:13", 4.04, "14/10/2023 14:51", 3.83, "13/10/2023 02:13", 4.86, "14/10/2023 14:51", "https://www.betexplorer.com/football/denmark/2nd-division/ab-copenhagen-skive/WKvgvmDm/")
Set-RowRange $ws 66 $row66
$row67 = @("Brabrand", 2, "AB Copenhagen", 2, 3.31, "20/10/2023 01:13", 3.55, "21/10/2023 13:41", 3.44, "20/10/2023 01:13", 3.49, "21/10/2023 13:41", 1.95, "20/10/2023 01:13", 1.99, "21/10/2023 13:41", "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-ab-copenhagen/Iys5yoc6/")
Set-RowRange $ws 67 $row67
$row69 = @("Middelfart", 2, "Nykobing", 2, 1.95, "20/10/2023 01:13", 2.1, "21/10/2023 13:41", 3.43, "20/10/2023 01:13", 3.66, "21/10/2023 13:41", 3.21, "20/10/2023 01:13", 3.12, "21/10/2023 13:41", "https://www.betexplorer.com/football/denmark/2nd-division/middelfart-nykobing/CzzEZORI/")
Set-RowRange $ws 69 $row69
$row70 = @("Skive", 1, "Roskilde", 2, 3.5, "20/10/2023 01:13", 3.25, "21/10/2023 13:41", 3.38, "20/10/2023 01:13", 3.43, "21/10/2023 13:41", 1.91, "20/10/2023 01:13", 2.12, "21/10/2023 13:41", "https://www.betexplorer.com/football/denmark/2nd-division/skive-roskilde/vor1xRs0/")
Set-RowRange $ws 70 $row70

# --- Append new rows 92-96 (copy formatting from row 91, then set values) ---
$ws.Range("A91:V91").Copy()
$ws.Range("A92:V96").PasteSpecial(-4122)

# Row 92
$ws.Range("A92").Value = 91
$ws.Range("B92").Value = "denmark"
$ws.Range("C92").Value = "2nd-division"
$ws.Range("D92").Value = "2023-2024"
$ws.Range("E92").Value = 45248.54166666666
$newrow92 = @("AB Copenhagen", 3, "FA 2000", 2, 1.71, "17/11/2023 01:11", 1.67, "18/11/2023 12:58", 3.78, "17/11/2023 01:11", 4.03, "18/11/2023 12:58", 3.75, "17/11/2023 01:11", 4.48, "18/11/2023 12:58", "https://www.betexplorer.com/football/denmark/2nd-division/ab-copenhagen-frederiksberg-alliancen-2000/faV42E1M/")
Set-RowRange $ws 92 $newrow92

# Row 93
$ws.Range("A93").Value = 92
$ws.Range("B93").Value = "denmark"
$ws.Range("C93").Value = "2nd-division"
$ws.Range("D93").Value = "2023-2024"
$ws.Range("E93").Value = 45248.54166666666
$newrow93 = @("Hellerup", 1, "Esbjerg", 4, 5.84, "17/11/2023 01:11", 7.11, "18/11/2023 12:57", 5.1, "17/11/2023 01:11", 5.63, "18/11/2023 12:59", 1.33, "17/11/2023 01:11", 1.32, "18/11/2023 12:57", "https://www.betexplorer.com/football/denmark/2nd-division/hellerup-esbjerg/QJwl6d2c/")
Set-RowRange $ws 93 $newrow93

# Row 94
$ws.Range("A94").Value = 93
$ws.Range("B94").Value = "denmark"
$ws.Range("C94").Value = "2nd-division"
$ws.Range("D94").Value = "2023-2024"
$ws.Range("E94").Value = 45248.58333333334
$newrow94 = @("Middelfart", 0, "Skive", 1, 1.57, "17/11/2023 02:12", 1.64, "18/11/2023 13:54", 3.9, "17/11/2023 02:12", 3.89, "18/11/2023 13:54", 4.7, "17/11/2023 02:12", 4.92, "18/11/2023 13:51", "https://www.betexplorer.com/football/denmark/2nd-division/middelfart-skive/MoXd4zX9/")
Set-RowRange $ws 94 $newrow94

# Row 95
$ws.Range("A95").Value = 94
$ws.Range("B95").Value = "denmark"
$ws.Range("C95").Value = "2nd-division"
$ws.Range("D95").Value = "2023-2024"
$ws.Range("E95").Value = 45248.58333333334
$newrow95 = @("Brabrand", 0, "Aarhus Fremad", 0, 4.09, "17/11/2023 02:12", 5.04, "18/11/2023 13:49", 3.99, "17/11/2023 02:12", 4.14, "18/11/2023 13:49", 1.61, "17/11/2023 02:12", 1.58, "18/11/2023 13:49", "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-aarhus-fremad/8rTh5GH3/")
Set-RowRange $ws 95 $newrow95

# Row 96
$ws.Range("A96").Value = 95
$ws.Range("B96").Value = "denmark"
$ws.Range("C96").Value = "2nd-division"
$ws.Range("D96").Value = "2023-2024"
$ws.Range("E96").Value = 45248.625
$newrow96 = @("Thisted FC", 1, "Roskilde", 3, 3.16, "17/11/2023 03:11", 3.28, "18/11/2023 14:56", 3.46, "17/11/2023 03:11", 3.57, "18/11/2023 14:56", 2, "17/11/2023 03:11", 2.05, "18/11/2023 14:56", "https://www.betexplorer.com/football/denmark/2nd-division/thisted-fc-roskilde/G8W03fnG/")
Set-RowRange $ws 96 $newrow96

